# Insert a new weekly price row at row 44 (pushing existing rows 44-71 down to 45-72)
# for the "Espinaca" sheet, adding the 2022-04-29 (serial 44680) observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 44 (Excel shifts rows 44:71 -> 45:72)
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new observation
$ws.Cells.Item(44, 1).Value  = 11
$ws.Cells.Item(44, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(44, 3).Value  = "Bíobío"
$ws.Cells.Item(44, 4).Value  = 44680
$ws.Cells.Item(44, 5).Value  = 8
$ws.Cells.Item(44, 6).Value  = 100112012
$ws.Cells.Item(44, 7).Value  = "Espinaca"
$ws.Cells.Item(44, 8).Value  = "Sin especificar"
$ws.Cells.Item(44, 9).Value  = "Primera"
$ws.Cells.Item(44, 10).Value = 50
$ws.Cells.Item(44, 11).Value = 6000
$ws.Cells.Item(44, 12).Value = 6500
$ws.Cells.Item(44, 13).Value = 6300
$ws.Cells.Item(44, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(44, 15).Value = "Región Metropolitana"
$ws.Cells.Item(44, 16).Value = 630
$ws.Cells.Item(44, 17).Value = 10
$ws.Cells.Item(44, 18).Value = "Hortaliza"
